$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 19 new rows (432-450) for the 2025-10-01 session; formatting (styles)
# is copied down automatically from row 431 (the prior last row) by Insert().
$ws.Range("A432:A450").EntireRow.Insert()

# Row 432: Ilan Ihaddadene
$ws.Range("A432").Value = 45931
$ws.Range("B432").Value = "Ilan Ihaddadene"
$ws.Range("C432").Value = 70
$ws.Range("D432").Value = 7
$ws.Range("E432").Value = 7
$ws.Range("F432").Value = 0
$ws.Range("H432").Value = 5
$ws.Range("I432").Formula = "=C432*D432"

# Row 433: Levy Ndoutoume
$ws.Range("A433").Value = 45931
$ws.Range("B433").Value = "Levy Ndoutoume"
$ws.Range("C433").Value = 70
$ws.Range("D433").Value = 7
$ws.Range("E433").Value = 7
$ws.Range("F433").Value = 5
$ws.Range("G433").Value = "Ischio"
$ws.Range("H433").Value = 4
$ws.Range("I433").Formula = "=C433*D433"

# Row 434: Hedi Nasri
$ws.Range("A434").Value = 45931
$ws.Range("B434").Value = "Hedi Nasri"
$ws.Range("C434").Value = 70
$ws.Range("D434").Value = 6
$ws.Range("E434").Value = 3
$ws.Range("F434").Value = 3
$ws.Range("G434").Value = "Adducteurs "
$ws.Range("H434").Value = 5
$ws.Range("I434").Formula = "=C434*D434"

# Row 435: Amine Taiar
$ws.Range("A435").Value = 45931
$ws.Range("B435").Value = "Amine Taiar"
$ws.Range("C435").Value = 70
$ws.Range("D435").Value = 4
$ws.Range("E435").Value = 6
$ws.Range("F435").Value = 6
$ws.Range("G435").Value = "Adducteur gauche"
$ws.Range("H435").Value = 2
$ws.Range("I435").Formula = "=C435*D435"

# Row 436: Jeremie Laurent
$ws.Range("A436").Value = 45931
$ws.Range("B436").Value = "Jeremie Laurent"
$ws.Range("C436").Value = 70
$ws.Range("D436").Value = 8
$ws.Range("E436").Value = 7
$ws.Range("F436").Value = 0
$ws.Range("H436").Value = 7
$ws.Range("I436").Formula = "=C436*D436"

# Row 437: Yoan Zouma
$ws.Range("A437").Value = 45931
$ws.Range("B437").Value = "Yoan Zouma"
$ws.Range("C437").Value = 70
$ws.Range("D437").Value = 4
$ws.Range("E437").Value = 7
$ws.Range("F437").Value = 5
$ws.Range("G437").Value = "Cheville"
$ws.Range("H437").Value = 5
$ws.Range("I437").Formula = "=C437*D437"

# Row 438: Yoann Martelat
$ws.Range("A438").Value = 45931
$ws.Range("B438").Value = "Yoann Martelat"
$ws.Range("C438").Value = 70
$ws.Range("D438").Value = 6
$ws.Range("E438").Value = 7
$ws.Range("F438").Value = 5
$ws.Range("G438").Value = "Genou"
$ws.Range("H438").Value = 5
$ws.Range("I438").Formula = "=C438*D438"

# Row 439: Amir Etien
$ws.Range("A439").Value = 45931
$ws.Range("B439").Value = "Amir Etien"
$ws.Range("C439").Value = 70
$ws.Range("D439").Value = 6
$ws.Range("E439").Value = 7
$ws.Range("F439").Value = 6
$ws.Range("G439").Value = "Genou "
$ws.Range("H439").Value = 3
$ws.Range("I439").Formula = "=C439*D439"

# Row 440: Ilyes Boughanmi
$ws.Range("A440").Value = 45931
$ws.Range("B440").Value = "Ilyes Boughanmi"
$ws.Range("C440").Value = 70
$ws.Range("D440").Value = 6
$ws.Range("E440").Value = 6
$ws.Range("F440").Value = 0
$ws.Range("H440").Value = 0
$ws.Range("I440").Formula = "=C440*D440"

# Row 441: Omar Benyounes
$ws.Range("A441").Value = 45931
$ws.Range("B441").Value = "Omar Benyounes"
$ws.Range("C441").Value = 70
$ws.Range("D441").Value = 4
$ws.Range("E441").Value = 5
$ws.Range("F441").Value = 0
$ws.Range("H441").Value = 3
$ws.Range("I441").Formula = "=C441*D441"

# Row 442: Kamal Bafounta
$ws.Range("A442").Value = 45931
$ws.Range("B442").Value = "Kamal Bafounta"
$ws.Range("C442").Value = 70
$ws.Range("D442").Value = 6
$ws.Range("E442").Value = 5
$ws.Range("F442").Value = 3
$ws.Range("G442").Value = "Cheville"
$ws.Range("H442").Value = 7
$ws.Range("I442").Formula = "=C442*D442"

# Row 443: Malik Boussaid
$ws.Range("A443").Value = 45931
$ws.Range("B443").Value = "Malik Boussaid"
$ws.Range("C443").Value = 70
$ws.Range("D443").Value = 3
$ws.Range("E443").Value = 1
$ws.Range("F443").Value = 0
$ws.Range("H443").Value = 10
$ws.Range("I443").Formula = "=C443*D443"

# Row 444: Naim Ighbane
$ws.Range("A444").Value = 45931
$ws.Range("B444").Value = "Naim Ighbane"
$ws.Range("C444").Value = 70
$ws.Range("D444").Value = 4
$ws.Range("E444").Value = 7
$ws.Range("F444").Value = 2
$ws.Range("G444").Value = "Cheville"
$ws.Range("H444").Value = 4
$ws.Range("I444").Formula = "=C444*D444"

# Row 445: Karim Belmahi
$ws.Range("A445").Value = 45931
$ws.Range("B445").Value = "Karim Belmahi"
$ws.Range("C445").Value = 70
$ws.Range("D445").Value = 7
$ws.Range("E445").Value = 8
$ws.Range("F445").Value = 0
$ws.Range("H445").Value = 10
$ws.Range("I445").Formula = "=C445*D445"

# Row 446: Emmanuel Valey
$ws.Range("A446").Value = 45931
$ws.Range("B446").Value = "Emmanuel Valey"
$ws.Range("C446").Value = 70
$ws.Range("D446").Value = 7
$ws.Range("E446").Value = 5
$ws.Range("F446").Value = 6
$ws.Range("G446").Value = "Cheville coup"
$ws.Range("H446").Value = 7
$ws.Range("I446").Formula = "=C446*D446"

# Row 447: Karahali Souaré
$ws.Range("A447").Value = 45931
$ws.Range("B447").Value = "Karahali Souaré"
$ws.Range("C447").Value = 70
$ws.Range("D447").Value = 3
$ws.Range("E447").Value = 5
$ws.Range("F447").Value = 7
$ws.Range("G447").Value = "Cheville "
$ws.Range("H447").Value = 4
$ws.Range("I447").Formula = "=C447*D447"

# Row 448: Sofiane Belle
$ws.Range("A448").Value = 45931
$ws.Range("B448").Value = "Sofiane Belle"
$ws.Range("C448").Value = 70
$ws.Range("D448").Value = 5
$ws.Range("E448").Value = 3
$ws.Range("F448").Value = 0
$ws.Range("H448").Value = 3
$ws.Range("I448").Formula = "=C448*D448"

# Row 449: Mattheo Haon
$ws.Range("A449").Value = 45931
$ws.Range("B449").Value = "Mattheo Haon"
$ws.Range("C449").Value = 70
$ws.Range("D449").Value = 8
$ws.Range("E449").Value = 8
$ws.Range("F449").Value = 0
$ws.Range("H449").Value = 7
$ws.Range("I449").Formula = "=C449*D449"

# Row 450: Naim Dhib
$ws.Range("A450").Value = 45931
$ws.Range("B450").Value = "Naim Dhib"
$ws.Range("C450").Value = 70
$ws.Range("D450").Value = 6
$ws.Range("E450").Value = 7
$ws.Range("F450").Value = 1
$ws.Range("G450").Value = "Jambes"
$ws.Range("H450").Value = 5
$ws.Range("I450").Formula = "=C450*D450"

# The "Localisation douleur" (G) cells above that contain text should use the
# same style as other populated cells in that column (copied from row 430,
# a nearby row whose G cell already holds text), rather than the blank-cell
# style that Insert() propagated from row 431 (whose G cell was empty).
$ws.Range("G430").Copy()
$ws.Range("G433").PasteSpecial(-4122)
$ws.Range("G434").PasteSpecial(-4122)
$ws.Range("G435").PasteSpecial(-4122)
$ws.Range("G437").PasteSpecial(-4122)
$ws.Range("G438").PasteSpecial(-4122)
$ws.Range("G439").PasteSpecial(-4122)
$ws.Range("G442").PasteSpecial(-4122)
$ws.Range("G444").PasteSpecial(-4122)
$ws.Range("G446").PasteSpecial(-4122)
$ws.Range("G447").PasteSpecial(-4122)
$ws.Range("G450").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the selection to match the latest data (the user had scrolled down
# and selected the cell below the newly-entered rows).
$ws.Range("K444").Select()
